$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: author name change
$ws.Range("A2").Value = "Héléna Reymond"

# Row 5 -> becomes the new "Recherche d'idées de projets en classe" entry (0.25h)
# Row 6 -> the old row5 content shifts here but with a new date (43151) and new text
# Row 7 -> the old row6 content shifts here but with a new date (43151) and new text

# Set row 6 (will hold the "fonctionnalités" text) height to match old row5 (wrap-computed 43.2)
$ws.Rows.Item(6).RowHeight = 43.2
$ws.Range("A6").Value = 43151
$ws.Range("B6").Value = "Recherche d'idée de fonctionnalités pour l'application de budget (analyse d'applications existantes, réflexions personnelles)"
$ws.Range("C6").Value = 0.5

# Set row 7 (will hold the "cahier des charges" text) height to match old row6 (wrap-computed 28.8)
$ws.Rows.Item(7).RowHeight = 28.8
$ws.Range("A7").Value = 43151
$ws.Range("B7").Value = "Co-rédaction, correction et relecture du cahier des charges réalisé par Daniel"
$ws.Range("C7").Value = 0.5

# Row 5 becomes the new short entry; height goes back to the default (no custom height)
$ws.Range("B5").Value = "Recherche d'idées de projets en classe"
$ws.Range("C5").Value = 0.25
$ws.Rows.Item(5).AutoFit()

# Selection moved to B6 in the saved file
[void]$ws.Range("B6").Select()
